$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 1 new row at row 142
$ws.Rows(142).Insert()
$ws.Range("B142").Value = 5170
$ws.Range("C142").Formula = "=IF(B142=`$E`$1,1,0)"

# Insert 2 new rows at row 224
$ws.Range("224:225").Insert()
$ws.Range("B224").Value = 5164
$ws.Range("C224").Formula = "=IF(B224=`$E`$1,1,0)"
$ws.Range("B225").Value = 2354
$ws.Range("C225").Formula = "=IF(B225=`$E`$1,1,0)"

# Row 224 joins the merged label group above it (A220:A223 -> A220:A224)
$ws.Range("A220:A224").Merge()

# Update view state
$ws.Application.ActiveWindow.ScrollRow = 107
$ws.Range("B142").Select()
